# =====================================================================
# CV update: rewritten summary, expanded bullet detail, $XXM placeholder
# cleanup, and de-duplicated PSO-practices / JAPAC-SEA scope language.
# =====================================================================
$d = $word.ActiveDocument

function Set-ParaTextByPrefix($Prefix, $NewText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like ($Prefix + "*")) {
            $p.Range.Text = $NewText
            return $i
        }
    }
    throw ("Paragraph not found for prefix: " + $Prefix)
}

function Insert-ParaAfterIndex($Index, $NewText) {
    $p = $d.Paragraphs.Item($Index)
    $p.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($Index + 1)
    $newP.Range.Text = $NewText
    return ($Index + 1)
}

# --- 1. Professional Summary: lead with Data & AI / team-building framing ---
Set-ParaTextByPrefix 'Senior Engineering Leader with 15+ years of experience bridg' 'Senior Engineering Leader with 15+ years building the teams, frameworks, and systems that turn Data and AI from research to production. Currently Head of Data & Analytics for Google Cloud in Southeast Asia - a practice built from zero, delivering enterprise Data and AI transformation across 7 countries.  Dual track as "Player-Coach": leading petabyte-scale data platforms and multi-agent systems for Fortune 500 clients, while driving innovation through published research (5 technical disclosures, 6 published packages on PyPI and Maven Central, plus open-source AI safety tools including sandbagging detection and activation steering). Member of Google Cloud delta, architecting solutions at the intersection of applied AI and enterprise scale.' | Out-Null

# --- 2. Head of D&A role blurb: JAPAC delivery scope / SEA site-lead scope ---
Set-ParaTextByPrefix 'Dual-track role combining technical innovation leadership wi' 'Dual-track role combining technical innovation leadership with regional delivery management. Built Google Cloud''s Data Analytics practice for Southeast Asia with delivery scope across JAPAC, while serving as Site Lead overseeing cross-practice operations in SEA. Member of <a href="https://cloud.google.com/consulting/innovation-and-transformation">delta</a> - Google Cloud''s innovation and transformation team architecting enterprise AI solutions at scale.' | Out-Null

# --- 3. Replace $XXM+ placeholders with descriptive language ---
Set-ParaTextByPrefix 'Direct $XXM+ Data Analytics delivery portfolio across JAPAC ' 'Direct regional Data Analytics delivery portfolio across JAPAC while simultaneously overseeing cross-practice portfolio as Site Lead.' | Out-Null

# --- 4. Drop redundant 7-PSO-practices parenthetical (already spelled out above) ---
Set-ParaTextByPrefix 'Pioneered agentic AI adoption across all 7 PSO practices (Da' 'Pioneered agentic AI adoption across all 7 PSO practices and 6 JAPAC sub-regions, building SDKs, agent catalog, and standardized templates while designing reusable governance frameworks that accelerated innovation and reduced delivery costs.' | Out-Null

# --- 5. Drop redundant "multi-million dollar" (already said of same engagements) ---
Set-ParaTextByPrefix 'Built agentic tool suites including architecture discovery (' 'Built agentic tool suites including architecture discovery (100M+ node graph modeling), automated data pipeline generation, and platform cleanup agents that recovered at-risk engagements and secured significant long-term cloud commitments.' | Out-Null

# --- 6. Standard Chartered bullets: rewrite intro + expand from 2 to 4 bullets ---
Set-ParaTextByPrefix 'Led enterprise-wide AI and data platform development serving' 'Led design and development of retail bank''s data & analytics platform serving 11 markets, 100+ systems, and 1200+ users.' | Out-Null
$idx = Set-ParaTextByPrefix 'Delivered a Self-Service ML Platform that reduced model deve' 'Developed self-service ML Workbench reducing model deployment time from months to weeks'
$idx = Set-ParaTextByPrefix 'MarTech modernization - +30% customer acquisition' 'Architected MarTech strategy driving 30% increase in customer acquisition through data-driven personalization'
$idx = Insert-ParaAfterIndex $idx 'Created credit risk models over 15,000+ named entities leveraging news trends and social signals, reducing potential losses by $5M'
$idx = Insert-ParaAfterIndex $idx 'Defined enterprise data strategy including third-party data governance, privacy frameworks, and cloud adoption roadmap'

# --- 7. Teradata bullets: expand from 2 to 4 bullets ---
$idx = Set-ParaTextByPrefix 'Data lakes processing 1.2 PB/hour for Fortune 500 clients ac' 'Designed 5 global data lakes with ETL pipelines handling 1.2 PB/hour and 40K daily files'
$idx = Set-ParaTextByPrefix 'Real-time fraud detection systems - 60% reduction in false p' 'Engineered real-time platform processing 2.5M events/second, improving Ad campaign responsiveness by 80%'
$idx = Insert-ParaAfterIndex $idx 'Built ML fraud detection system achieving 60% fewer false positives and 25% higher detection rates, resulting in $3M savings'
$idx = Insert-ParaAfterIndex $idx 'Built and managed large-scale Hadoop clusters (300+ nodes) for banks and telcos across JAPAC'

# --- 8. AI Metacognition Toolkit: swap "steering vectors" detail for PyPI credit ---
Set-ParaTextByPrefix 'Activation-level detection of sandbagging, deception, and si' 'Activation-level detection of sandbagging, deception, and situational awareness in LLMs. Linear probes achieve 90-96% accuracy across Mistral, Gemma, and Qwen models. Published on PyPI.' | Out-Null

Write-Host "CV updated."
